$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.408.21"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "2.408.52"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  +0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.08"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.95"
$ws.Range("E6").Value = "  +2.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("E8").Value = "  +1.48%  "
$ws.Range("D9").Value = "2.419.75"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("E10").Value = "  +2.30%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.23"
$ws.Range("E12").Value = "  +3.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.347"
$ws.Range("E13").Value = "  +3.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.59"
$ws.Range("E14").Value = "  +2.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000173"
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("D16").Value = "2.838.57"
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").Value = "61.287.80"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "2.422.06"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.12"
$ws.Range("E19").Value = "  +3.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.75"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.41"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.96"
$ws.Range("E24").Value = "  +7.07%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.31"
$ws.Range("E26").Value = "  +1.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "618.46"
$ws.Range("E27").Value = "  +6.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.30"
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("D29").Value = "0.0₃0954"
$ws.Range("E29").Value = "  +3.30%  "
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.03"
$ws.Range("E31").Value = "  +2.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.39"
$ws.Range("E32").Value = "  +3.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.81"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.133"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.48"
$ws.Range("E35").Value = "  +4.62%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.374"
$ws.Range("E37").Value = "  +2.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.63"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "151.31"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.33"
$ws.Range("E40").Value = "  +4.05%  "
$ws.Range("B41").Value = "EthereumClassic"
$ws.Range("C41").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.42"
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.71"
$ws.Range("E42").Value = "  +3.11%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.54"
$ws.Range("E43").Value = "  +6.09%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.95"
$ws.Range("E45").Value = "  +1.83%  "
$ws.Range("D46").Value = "0.0₆0286"
$ws.Range("E46").Value = "  -2.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.00"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.54"
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.95"
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.594"
$ws.Range("E50").Value = "  +1.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0512"
$ws.Range("E51").Value = "  +2.68%  "
